# Update view-count figures in column F that changed between scrapes.
# Sheet "展览" (exhibitions): F5, F6, F12, F14
# Sheet "全部类型" (all types, duplicates the same rows): F6, F7, F13, F15

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 6995
$wsExpo.Range("F6").Value = 2512
$wsExpo.Range("F12").Value = 37
$wsExpo.Range("F14").Value = 579

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 6995
$wsAll.Range("F7").Value = 2512
$wsAll.Range("F13").Value = 37
$wsAll.Range("F15").Value = 579
